# The spreadsheet tracks a SKOS vocabulary. The Google-sheet regenerated the
# .ttl export: the version label was simplified, a stray "TestTerm" row was
# removed (so every concept below it shifts up one row), and a handful of
# skos:prefLabel values for the rows that moved picked up a space between
# the words of their CamelCase identifiers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) pav:version value changes from "version 0.1" to "0.1.0"
$ws.Range("B13").Value = "0.1.0"

# 2) Remove the obsolete "rock-n-roll:TestTerm" row (row 18).
#    This shifts every following row up by one, which is exactly what the
#    diff shows happening to rows 19-28 (now 18-27), including the removal
#    of what used to be the trailing empty row 28.
$ws.Rows.Item(18).Delete()

# 3) A few of the rows that shifted up also got their skos:prefLabel
#    (column B) reformatted from CamelCase into space-separated words.
$ws.Range("B19").Value = "Geologic Concept"
$ws.Range("B20").Value = "Geologic Property"
$ws.Range("B21").Value = "Earth Material"
$ws.Range("B22").Value = "Compound Material"
